# "writeOff box for ENG"
# - add a new "tubes" sheet (right after "ready_to_sale", before "components")
# - update stock counters on "ready_to_sale" (ENG) and "components"
# - drop two obsolete rows ("Пенка 130 см" / "Пенка 250 см") from "components"
# - minor header line-break fix on "Sales"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) insert the new "tubes" sheet right before "components" FIRST -
#    sheet references are positional, so grab fresh handles for every
#    other sheet only *after* the sheet collection has been resized.
# ---------------------------------------------------------------------
$wsComponentsBeforeInsert = $wb.Worksheets.Item("components")
$wsTubes = $wb.Worksheets.Add($wsComponentsBeforeInsert)
$wsTubes.Name = "tubes"

$tubeNames = @("Fire", "Waterfall", "Ether-Acril", "Ether-Wood", "Alchemy", "Infinity", "Eternal-love")
$r = 2
foreach ($tubeName in $tubeNames) {
    $wsTubes.Range("A$r").Value = $tubeName
    $r = $r + 1
}

# re-fetch all other sheets by name now that the insert has happened
$wsReady = $wb.Worksheets.Item("ready_to_sale")
$wsComponents = $wb.Worksheets.Item("components")
$wsSales = $wb.Worksheets.Item("Sales")

# ---------------------------------------------------------------------
# 2) ready_to_sale : write-off / restock adjustments
# ---------------------------------------------------------------------
$wsReady.Range("B2").Value = 7    # Fire        : В наличии ENG 5 -> 7
$wsReady.Range("C2").Value = 6    # Fire        : В наличии UA  0 -> 6
$wsReady.Range("B3").Value = 4    # Waterfall   : В наличии ENG 1 -> 4
$wsReady.Range("B6").Value = 2    # Alchemy     : В наличии ENG 4 -> 2
$wsReady.Range("C6").Value = 2    # Alchemy     : В наличии UA  5 -> 2
$wsReady.Range("B7").Value = 5    # Infinity    : В наличии ENG 3 -> 5
$wsReady.Range("C7").Value = 3    # Infinity    : В наличии UA  10 -> 3

# ---------------------------------------------------------------------
# 3) components : stock adjustments + drop two obsolete rows
# ---------------------------------------------------------------------
$wsComponents.Range("B2").Value = 23   # Bag стандарт      : 17 -> 23
$wsComponents.Range("B4").Value = 106  # Box Divya         : 110 -> 106
$wsComponents.Range("B5").Value = 23   # Планки дерево Б   : 17 -> 23
$wsComponents.Range("B6").Value = 23   # Планки дерево М   : 17 -> 23
$wsComponents.Range("B9").Value = 23   # Подставки         : 17 -> 23
$wsComponents.Range("B11").Value = 18  # Стики             : 12 -> 18
$wsComponents.Range("B13").Value = 10  # Флизелин Ефир     : 0 -> 10
$wsComponents.Range("B16").Value = 0   # Миникорд 110 см   : 5 -> 0

$wsComponents.Range("A20:B21").Delete()

# ---------------------------------------------------------------------
# 4) Sales : header wraps gained an extra line break
# ---------------------------------------------------------------------
$wsSales.Range("E2").Value = "Имя`n`n`nпокупателя"
$wsSales.Range("G2").Value = "Откуда`n`n`nпокупатель"
$wsSales.Range("H2").Value = "Track`n`n`nnumber"
$wsSales.Range("I2").Value = "Стоимость`n`n`nдоставки"
$wsSales.Range("J2").Value = "Почтовая`n`n`nслужба"

Write-Host "edit complete"
